$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D value swaps (re-imputed vs un-imputed cells) ---
$ws.Range("D2").Value = -13.5
$ws.Range("D3").ClearContents()
$ws.Range("D4").ClearContents()
$ws.Range("D11").Value = -15.5
$ws.Range("D13").ClearContents()
$ws.Range("D21").Value = -14.3
$ws.Range("D25").ClearContents()

# --- Remove row "RM 232" (row 26) entirely ---
$ws.Rows.Item(26).Delete()

# --- Remove row "SC 92" (now row 27 after the previous delete) ---
$ws.Rows.Item(27).Delete()

# --- Fill in the previously-missing B/D values for the final row "SC 232" ---
$ws.Range("B33").Value = -19.5
$ws.Range("D33").Value = -14.1

Write-Output "done"
